# Apply updated transition-probability values to the matrix sheet.
# (added more games, sped up simulate game logic, and drafted optimization logic)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.1428571428571428
    "C2" = 0.7142857142857143
    "S2" = 0.1428571428571428

    "C3" = 0.1428571428571428
    "P3" = 0.7142857142857143
    "S3" = 0.1428571428571428

    "S4" = 1

    "F6" = 0.08333333333333333
    "J6" = 0.2083333333333333
    "O6" = 0.04166666666666666
    "Q6" = 0.04166666666666666
    "R6" = 0.25
    "S6" = 0.375

    "B7" = 0.1111111111111111
    "D7" = 0.1111111111111111
    "O7" = 0.1111111111111111
    "Q7" = 0.3333333333333333
    "R7" = 0.1111111111111111
    "S7" = 0.2222222222222222

    "B8" = 0.0410958904109589
    "D8" = 0.0136986301369863
    "F8" = 0.0273972602739726
    "J8" = 0.1506849315068493
    "Q8" = 0.0958904109589041
    "R8" = 0.0684931506849315
    "S8" = 0.6027397260273972

    "B9" = 0.05555555555555555
    "J9" = 0.2222222222222222
    "R9" = 0.1666666666666667
    "S9" = 0.5555555555555556

    "B10" = 0.01265822784810127
    "D10" = 0.01265822784810127
    "F10" = 0.1392405063291139
    "J10" = 0.08860759493670886
    "Q10" = 0.1392405063291139
    "R10" = 0.0759493670886076
    "S10" = 0.5316455696202531

    "L11" = 1

    "G12" = 0.5454545454545454
    "J12" = 0.1818181818181818
    "L12" = 0.09090909090909091
    "S12" = 0.1818181818181818

    "G13" = 1

    "F15" = 0.07142857142857142
    "H15" = 0.07142857142857142
    "I15" = 0.1428571428571428
    "J15" = 0.4285714285714285
    "O15" = 0.07142857142857142
    "S15" = 0.2142857142857143

    "H16" = 0.25
    "J16" = 0.75

    "H17" = 0.4090909090909091
    "I17" = 0.09090909090909091
    "J17" = 0.2727272727272727
    "O17" = 0.04545454545454546
    "S17" = 0.1818181818181818

    "H18" = 0.35
    "I18" = 0.25
    "J18" = 0.15
    "K18" = 0.1
    "S18" = 0.15

    "F19" = 0.02702702702702703
    "H19" = 0.3783783783783784
    "I19" = 0.06081081081081081
    "J19" = 0.2297297297297297
    "K19" = 0.05405405405405406
    "M19" = 0.02027027027027027
    "O19" = 0.06081081081081081
    "S19" = 0.1689189189189189
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
